# OneFifthBox1DisplacementCals.xlsx - "New Shunt Cal sheets for the small scale LBCBs"
#
# The "Displacement Cals" sheet keeps three independent calibration blocks
# (rows 2-22, 23-43, 44-64). Each block computes a fresh regression slope
# ("New Slope:", row 17/38/59) and compares it against the slope measured
# during the previous calibration run ("Previous Slopes:", row 20/41/62 -
# these are hand-entered numbers, not formulas). The Delta/%Difference rows
# underneath are formulas and recompute automatically once the "Previous
# Slopes" numbers are refreshed.
#
# This edit refreshes the three "Previous Slopes:" rows with the values
# recorded from the latest calibration pass, and updates the sheet's saved
# selection to reflect where the user left off (cell F65, just below the
# data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Displacement Cals")

# --- Block 1 (rows 2-22) : X1 / X2 previous slopes -------------------------
$ws.Range("C20").Value = -0.26734361718872979
$ws.Range("D20").Value = 0.26732885080511981
$ws.Range("F20").Value = -0.27068181818181819
$ws.Range("G20").Value = 0.27179659542604756

# --- Block 2 (rows 23-43) : Y1 previous slopes ------------------------------
$ws.Range("C41").Value = -0.12935563380281689
$ws.Range("D41").Value = 0.12925345471864672
$ws.Range("F41").Value = -0.1312823660714286
$ws.Range("G41").Value = 0.13208895570119597

# --- Block 3 (rows 44-64) : Z1 / Z2 / Z3 previous slopes --------------------
$ws.Range("C62").Value = -0.13711495535714288
$ws.Range("D62").Value = 0.13783328328485747
$ws.Range("F62").Value = -0.13732254464285715
$ws.Range("G62").Value = 0.13748322978580291

# Leave the sheet scrolled/selected where the editor last left it.
$ws.Activate()
$ws.Range("F65").Select()
